$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the bug that was causing PayerCityStateZip (B9) to duplicate the
# PayeeCityStateZip value. Give it its own distinct value.
$ws.Range("B9").Value = "Austin, TX 78759"

# Update the active selection to reflect where the edit was made.
$ws.Activate()
$ws.Range("B9").Select()
